$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model names for A2:A26 (a re-shuffled/renumbered ordering of model_28_3_0..24)
$names = @(
    "model_28_3_0",
    "model_28_3_22",
    "model_28_3_21",
    "model_28_3_20",
    "model_28_3_19",
    "model_28_3_18",
    "model_28_3_17",
    "model_28_3_16",
    "model_28_3_15",
    "model_28_3_14",
    "model_28_3_13",
    "model_28_3_23",
    "model_28_3_12",
    "model_28_3_10",
    "model_28_3_9",
    "model_28_3_8",
    "model_28_3_7",
    "model_28_3_6",
    "model_28_3_5",
    "model_28_3_4",
    "model_28_3_3",
    "model_28_3_2",
    "model_28_3_1",
    "model_28_3_11",
    "model_28_3_24"
)

# New metric values (B..Q), identical across every row after the update
$metrics = @(
    0.9999549410293126,
    0.9990639202211722,
    0.9999877189213422,
    0.9999999999995537,
    0.9999924619841859,
    0.00004206057127643839,
    0.0008737893843822983,
    0.000009903151542305334,
    0.0000000000001813269564866554,
    0.000004951575861816146,
    0.0004236862227404525,
    0.006485412190172526,
    1.0000514959665,
    0.00676150961376217,
    110.1527996107309,
    165.0022117297999
)

for ($i = 0; $i -lt 25; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($j = 0; $j -lt 16; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $metrics[$j]
    }
}
